$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "In order to understand more clearly results of the survey, we decided to find a way to analyze them, comparing what person which ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "In order to understand more clearly the results of the survey, we decided to analyze them by looking for relationships in the data, such as how people who answered one question in a certain way answered other questions.",
    2)

$d.Content.Find.Execute(
    "answered specifically in one question said in the rest of question. That could help us finding the main issue with licences for people who care about the legality.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "This could help us find the main issue and define the problem.",
    2)

$d.Content.Find.Execute(
    "Due to the significant number of results we will not be able to do that manually. We can try to create program which will do that,",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Due to the significant number of results we will not be able to co-relate answers manually. We can try to create a program which will do that,",
    2)

$d.Content.Find.Execute(
    "however we may lack of time before the mid-term presentation. If we manage to do that, it may turn out to be resourceful",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "however we may lack time before the mid-term presentation.",
    2)

$d.Content.Find.Execute(
    "as a team, what are possible presentation structure, how to maintain our presentation in Pecha-Kucha form, what parts of our work should be emphasised and many more.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "as a team, what are the possible presentation structures, how to maintain our presentation in Pecha-Kucha form, what parts of our work should be emphasised and many more.",
    2)

$d.Content.Find.Execute(
    "and how to prepare so well, that everything will go smoothly. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "and how to prepare well, so that everything will go smoothly. ",
    2)

$d.Content.Find.Execute(
    "We can do everything in order to reduce our stress. There are really no rules about the way we should split the talking between ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "We can do multiple things in order to reduce our stress. There are really no rules about the way we should split the talking between ",
    2)

$d.Content.Find.Execute(
    "by the projector. We should made the IDs for each one of us.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "by the projector. We should make IDs for each one of us.",
    2)
